$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hashcode.csv")

$ws.Range("B9").Value = "6c1a51b5e85289c72b553ad899db1a51"
$ws.Range("B11").Value = "1f682c4baf00039722b9d3b2a8f6431f"
$ws.Range("B15").Value = "748fdfa86f292b105e0f90f30045b1c7"
$ws.Range("B24").Value = "d6e0c50a94994e93363908ad3893b5fb"
$ws.Range("B29").Value = "14bcaddadb80968b9f4d54b9a0bf5b4d"
$ws.Range("B34").Value = "9b5fa738b68a8c46f512c3e8ae609d3b"
$ws.Range("B121").Value = "81667d4f5140992663fc6287a415e11f"
$ws.Range("B133").Value = "e67cb7acf6fa4ff9ebae00920bc5988a"
$ws.Range("B136").Value = "4f5900008902af644f9584451a1c3901"
$ws.Range("B162").Value = "537a5222143850acb0b8e7c2a56d1a6f"
$ws.Range("B175").Value = "1aab84fff0421b283487ec2cd09a72ec"
$ws.Range("B180").Value = "8e3e66726412138b9c21d57bc4009d98"
$ws.Range("B183").Value = "b2ea477540860dd093eec216119c4402"
$ws.Range("B191").Value = "6aeb8c7ff9726e431785864e912f5be0"
$ws.Range("B198").Value = "d3bf3c1c93e8e11b73485fcb6846fff5"
$ws.Range("B213").Value = "f1a3da6a4991d211f4d0e18b9486ed7a"
$ws.Range("B228").Value = "da137e8bd5d8f137f86514581a664b40"
$ws.Range("B232").Value = "ff4d4a34187f39c569ebefb41b8d6aa3"
$ws.Range("B461").Value = "060072cb4a449d58d07838c00b609f70"
$ws.Range("B480").Value = "c2cefcf8311326ea2d84c3e9ddd5d4ad"
$ws.Range("B501").Value = "b960af1343abd623aaf02c982f837560"
$ws.Range("B506").Value = "aa1791820592e49d2dde3aff5748084a"
$ws.Range("B514").Value = "0163ad4ebad868ebcb1fb1d515410e6b"
$ws.Range("B524").Value = "b8463e643f40c14c051b7aa3e19cc647"
$ws.Range("B563").Value = "58aeeda8ebd6873d630280821cb636b9"
$ws.Range("B572").Value = "5ed55f8b2ae0bd9cea467720286f267b"
$ws.Range("B629").Value = "b4bf40be839e72ff90e5a588136c4567"
$ws.Range("B649").Value = "ea19ac78d9def67994c85b8b5c27e9a8"
$ws.Range("B655").Value = "e5f700c8b43c086d0c838f66e1305e35"
$ws.Range("B666").Value = "d0198b482e7ad0701fea272aba6657a8"
$ws.Range("B680").Value = "902b8c6b60528c7b830052fb906e50a8"
$ws.Range("B685").Value = "225498260d678337a4782766a1156652"
$ws.Range("B700").Value = "54f8f0d13d2be919db718fbd6002f7de"
$ws.Range("B703").Value = "09d87cbc478370a8a2f110e3b1786283"
$ws.Range("B704").Value = "0c15d29fc30a8a3b76d70a057ca66b27"
$ws.Range("B715").Value = "6fec891a7daf86028b2467a7fac67a3f"
$ws.Range("B729").Value = "b4db0bd5cfe9f51ea71702c7935a8b82"
$ws.Range("B742").Value = "3945cc1ced32bc3ccd9b183feb1b5bcd"
$ws.Range("B830").Value = "f1a61ae09a06993f94701cb2daa2fa3d"
$ws.Range("B835").Value = "493485141f8ff34952434469c68d6932"
$ws.Range("B854").Value = "36b9c870f13364b4b979e67b80e9bd2b"
$ws.Range("B862").Value = "27ca3467c0df769fe74125c62a70f180"
